$d = $word.ActiveDocument

# 1. Highlight (lightGray) the "Aggiungi nel file configurazione..." paragraph,
#    which spans two runs: the main sentence and "Edi.java".
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Aggiungi nel file configurazione l’intervallo con cui  i sensori inviano i dati di consumo e sistema di conseguenza Edi.java",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find 'Aggiungi nel file configurazione...' text" }
$r1.HighlightColorIndex = 16

# 2. Highlight (lightGray) the "Controlla gli intervalli..." paragraph, which
#    spans five runs (incl. the misspelled "poì" and "int").
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "Controlla gli intervalli di sicurezza e quelli di invio sensore: forse hai fatto un poì di confusione tra int e long",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find 'Controlla gli intervalli...' text" }
$r2.HighlightColorIndex = 16

# 3. Remove the whole "Crea l’interfaccia di usercmd e sistema di conseguenza il
#    file Edi" paragraph (now done, per the commit message).
$r3 = $d.Content
$found3 = $r3.Find.Execute(
    "Crea l’interfaccia di usercmd e sistema di conseguenza il file Edi",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find 'Crea l’interfaccia di usercmd...' text" }
# Re-fetch the containing paragraph from the document's Paragraphs collection so
# its Range includes the trailing paragraph mark, and deleting it merges the
# paragraph away entirely instead of leaving an empty one behind.
$paraIndex = $r3.Paragraphs.Item(1).Index
$d.Paragraphs.Item($paraIndex).Range.Delete()
